$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.703.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.84%  '
$ws.Range("D3").Value = "'1.874.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.36%  '
$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = "'282.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").Value = "'0.9989"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = "'0.5159"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.08%  '
$ws.Range("D8").Value = "'0.3531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.06%  '
$ws.Range("D9").Value = "'45.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.80%  '
$ws.Range("D10").Value = "'0.07168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.25%  '
$ws.Range("D11").Value = "'20.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = "'0.8229"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("D13").Value = "'0.07747"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = "'1.863.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.79%  '
$ws.Range("D15").Value = "'5.151"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.06%  '
$ws.Range("D16").Value = "'89.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("D17").Value = "'0.9987"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = "'14.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.89%  '
$ws.Range("D19").Value = "'0.000008205"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.88%  '
$ws.Range("D20").Value = "'0.9980"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").Value = "'26.761.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.73%  '
$ws.Range("D22").Value = "'4.805"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").Value = "'10.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D24").Value = "'6.252"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.93%  '
$ws.Range("D25").Value = "'2.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.79%  '
$ws.Range("D26").Value = "'145.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.21%  '
$ws.Range("D27").Value = "'17.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.52%  '
$ws.Range("D28").Value = "'1.668"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("E29").Value = '  +2.26%  '
$ws.Range("D30").Value = "'4.415"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("D31").Value = "'4.365"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.87%  '
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").Value = "'0.04912"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.51%  '
$ws.Range("D34").Value = "'1.179"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.23%  '
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = "'2.864"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = "'3.291"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.68%  '
$ws.Range("D38").Value = "'2.430"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("D39").Value = "'0.5301"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").Value = "'0.01880"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").Value = "'0.9735"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").Value = "'117.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.91%  '
$ws.Range("D43").Value = "'6.307"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.23%  '
$ws.Range("D44").Value = "'8.207"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").Value = "'0.4619"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("E47").Value = '  -1.49%  '
$ws.Range("D48").Value = "'9.462"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.18%  '
$ws.Range("D49").Value = "'36.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.78%  '
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").Value = "'0.05927"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.95%  '
